$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 349.22223
$ws.Range("I4").Value = 163.28572
$ws.Range("K4").Value = 163.28572
$ws.Range("M4").Value = -49.28572
$ws.Range("H28").Value = 350.04544
$ws.Range("I28").Value = 346.5625
$ws.Range("J28").Value = 359.33334
$ws.Range("K28").Value = 346.5625
$ws.Range("L28").Value = 359.33334
$ws.Range("M28").Value = 138.4375
$ws.Range("N28").Value = -1329.33334
$ws.Range("H31").Value = 1397.5
$ws.Range("I31").Value = 945
$ws.Range("K31").Value = 2835
$ws.Range("M31").Value = -2605
$ws.Range("H32").Value = 360.4
$ws.Range("J32").Value = 751
$ws.Range("L32").Value = 751
$ws.Range("N32").Value = -1403
$ws.Range("H98").Value = 2336.9473
$ws.Range("I98").Value = 1328.0769
$ws.Range("K98").Value = 1328.0769
$ws.Range("M98").Value = 169.9231
$ws.Range("H113").Value = 4385.8823
$ws.Range("I113").Value = 4074
$ws.Range("K113").Value = 4074
$ws.Range("M113").Value = -820
$ws.Range("H122").Value = 2336.9473
$ws.Range("I122").Value = 1328.0769
$ws.Range("K122").Value = 3984.2307
$ws.Range("M122").Value = -1534.2307
$ws.Range("H132").Value = 2440733.2
$ws.Range("I132").Value = 2704215.2
$ws.Range("J132").Value = 3524.75
$ws.Range("K132").Value = 8112645.600000001
$ws.Range("L132").Value = 10574.25
$ws.Range("M132").Value = -8110115.600000001
$ws.Range("N132").Value = -15634.25
$ws.Range("H137").Value = 1820459
$ws.Range("I137").Value = 2224432.5
$ws.Range("J137").Value = 2579
$ws.Range("K137").Value = 6673297.5
$ws.Range("L137").Value = 7737
$ws.Range("M137").Value = -6670747.5
$ws.Range("N137").Value = -12837

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1325.31
$ws.Range("I32").Value = 1197.1086
$ws.Range("J32").Value = 2799.625
$ws.Range("K32").Value = 1197.1086
$ws.Range("L32").Value = 2799.625
$ws.Range("M32").Value = -910.1086
$ws.Range("N32").Value = -3373.625
$ws.Range("H45").Value = 1105.68
$ws.Range("I45").Value = 1095.5
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 1095.5
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -718.5
$ws.Range("N45").Value = -2104
$ws.Range("H61").Value = 1544.6735
$ws.Range("I61").Value = 626.561
$ws.Range("J61").Value = 6250
$ws.Range("K61").Value = 626.561
$ws.Range("L61").Value = 6250
$ws.Range("M61").Value = -414.561
$ws.Range("N61").Value = -6674
$ws.Range("H74").Value = 615.37036
$ws.Range("I74").Value = 615.37036
$ws.Range("K74").Value = 615.37036
$ws.Range("M74").Value = 258.62964
$ws.Range("H77").Value = 615.37036
$ws.Range("I77").Value = 615.37036
$ws.Range("K77").Value = 3076.8518
$ws.Range("M77").Value = 1291.1482
$ws.Range("H122").Value = 3577.182
$ws.Range("I122").Value = 2891.5
$ws.Range("K122").Value = 8674.5
$ws.Range("M122").Value = -6224.5
$ws.Range("H132").Value = 1561.9839
$ws.Range("I132").Value = 1170.0209
$ws.Range("J132").Value = 2905.8572
$ws.Range("K132").Value = 3510.0627
$ws.Range("L132").Value = 8717.571599999999
$ws.Range("M132").Value = -980.0626999999999
$ws.Range("N132").Value = -13777.5716
$ws.Range("H136").Value = 1544.6735
$ws.Range("I136").Value = 626.561
$ws.Range("J136").Value = 6250
$ws.Range("K136").Value = 1879.683
$ws.Range("L136").Value = 18750
$ws.Range("M136").Value = 670.317
$ws.Range("N136").Value = -23850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1008.25
$ws.Range("I94").Value = 935.9091
$ws.Range("J94").Value = 1096.6666
$ws.Range("K94").Value = 935.9091
$ws.Range("L94").Value = 1096.6666
$ws.Range("M94").Value = -484.9091
$ws.Range("N94").Value = -1998.6666
$ws.Range("H134").Value = 1649.1857
$ws.Range("I134").Value = 1197.0984
$ws.Range("K134").Value = 3591.2952
$ws.Range("M134").Value = -1056.2952

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 51.375
$ws.Range("I7").Value = 34.2
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 34.2
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = 78.8
$ws.Range("N7").Value = -306
$ws.Range("H31").Value = 1925616
$ws.Range("I31").Value = 2779200.2
$ws.Range("J31").Value = 5051.3125
$ws.Range("K31").Value = 2779200.2
$ws.Range("L31").Value = 5051.3125
$ws.Range("M31").Value = -2778905.2
$ws.Range("N31").Value = -5641.3125
$ws.Range("H34").Value = 1925616
$ws.Range("I34").Value = 2779200.2
$ws.Range("J34").Value = 5051.3125
$ws.Range("K34").Value = 2779200.2
$ws.Range("L34").Value = 5051.3125
$ws.Range("M34").Value = -2778998.2
$ws.Range("N34").Value = -5455.3125
$ws.Range("H35").Value = 1164.4445
$ws.Range("I35").Value = 1164.4445
$ws.Range("K35").Value = 1164.4445
$ws.Range("M35").Value = -870.4445000000001
$ws.Range("H58").Value = 7354798
$ws.Range("I58").Value = 921.4259
$ws.Range("J58").Value = 35719750
$ws.Range("K58").Value = 921.4259
$ws.Range("L58").Value = 35719750
$ws.Range("M58").Value = -718.4259
$ws.Range("N58").Value = -35720156
$ws.Range("H132").Value = 1579.42
$ws.Range("I132").Value = 1195.3658
$ws.Range("J132").Value = 3329
$ws.Range("K132").Value = 3586.0974
$ws.Range("L132").Value = 9987
$ws.Range("M132").Value = -1056.0974
$ws.Range("N132").Value = -15047
$ws.Range("H134").Value = 1425.3556
$ws.Range("I134").Value = 766.5263
$ws.Range("J134").Value = 5001.857
$ws.Range("K134").Value = 2299.5789
$ws.Range("L134").Value = 15005.571
$ws.Range("M134").Value = 235.4211
$ws.Range("N134").Value = -20075.571
$ws.Range("H136").Value = 7354798
$ws.Range("I136").Value = 921.4259
$ws.Range("J136").Value = 35719750
$ws.Range("K136").Value = 2764.2777
$ws.Range("L136").Value = 107159250
$ws.Range("M136").Value = -214.2776999999996
$ws.Range("N136").Value = -107164350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1661.8438
$ws.Range("I131").Value = 1793.0625
$ws.Range("J131").Value = 1530.625
$ws.Range("K131").Value = 5379.1875
$ws.Range("L131").Value = 4591.875
$ws.Range("M131").Value = -339.1875
$ws.Range("N131").Value = -14671.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 37373.285
$ws.Range("I22").Value = 4888.5
$ws.Range("J22").Value = 50367.2
$ws.Range("K22").Value = 4888.5
$ws.Range("L22").Value = 50367.2
$ws.Range("M22").Value = -4359.5
$ws.Range("N22").Value = -51425.2
$ws.Range("H122").Value = 4122.3125
$ws.Range("I122").Value = 3022.5
$ws.Range("J122").Value = 4977.722
$ws.Range("K122").Value = 9067.5
$ws.Range("L122").Value = 14933.166
$ws.Range("M122").Value = -6617.5
$ws.Range("N122").Value = -19833.166
$ws.Range("H132").Value = 1919.0952
$ws.Range("I132").Value = 1503.5918
$ws.Range("J132").Value = 3373.3572
$ws.Range("K132").Value = 4510.7754
$ws.Range("L132").Value = 10120.0716
$ws.Range("M132").Value = -1980.7754
$ws.Range("N132").Value = -15180.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 29700
$ws.Range("J36").Value = 29700
$ws.Range("L36").Value = 29700
$ws.Range("N36").Value = -30824
$ws.Range("H122").Value = 3158.4443
$ws.Range("I122").Value = 2789.4736
$ws.Range("J122").Value = 4034.75
$ws.Range("K122").Value = 8368.4208
$ws.Range("L122").Value = 12104.25
$ws.Range("M122").Value = -5918.4208
$ws.Range("N122").Value = -17004.25
$ws.Range("H132").Value = 1801.9149
$ws.Range("I132").Value = 1181.5151
$ws.Range("J132").Value = 3264.2856
$ws.Range("K132").Value = 3544.5453
$ws.Range("L132").Value = 9792.856800000001
$ws.Range("M132").Value = -1014.5453
$ws.Range("N132").Value = -14852.8568
$ws.Range("H136").Value = 1853289
$ws.Range("I136").Value = 2174956.8
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 6524870.399999999
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = -6522320.399999999
$ws.Range("N136").Value = -16200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 669332.75
$ws.Range("I122").Value = 771291.7
$ws.Range("J122").Value = 6599.5
$ws.Range("K122").Value = 2313875.1
$ws.Range("L122").Value = 19798.5
$ws.Range("M122").Value = -2311425.1
$ws.Range("N122").Value = -24698.5
$ws.Range("H132").Value = 168664.77
$ws.Range("I132").Value = 224262.8
$ws.Range("J132").Value = 29669.666
$ws.Range("K132").Value = 672788.3999999999
$ws.Range("L132").Value = 89008.99800000001
$ws.Range("M132").Value = -670258.3999999999
$ws.Range("N132").Value = -94068.99800000001
$ws.Range("H136").Value = 1165.0667
$ws.Range("I136").Value = 836.03705
$ws.Range("J136").Value = 1658.6111
$ws.Range("K136").Value = 2508.11115
$ws.Range("L136").Value = 4975.8333
$ws.Range("M136").Value = 41.88884999999982
$ws.Range("N136").Value = -10075.8333
